$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Fill Sheet2 with the same card-data layout as Sheet1 -----------------
# Header row is identical to Sheet1's header.
$ws2.Range("A1").Value = $ws1.Range("A1").Value2
$ws2.Range("B1").Value = $ws1.Range("B1").Value2
$ws2.Range("C1").Value = $ws1.Range("C1").Value2
$ws2.Range("D1").Value = $ws1.Range("D1").Value2
$ws2.Range("E1").Value = $ws1.Range("E1").Value2

# Data rows: same numbers/names as Sheet1, but column A holds the card
# number as TEXT (e.g. "2-111111") instead of a plain number.
$cardNos = @("2-111111", "2-222222", "2-33333", "2-44444")
for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2
    $ws2.Range("A$r").Value = $cardNos[$i]
    $ws2.Range("B$r").Value = $ws1.Range("B$r").Value2
    $ws2.Range("C$r").Value = $ws1.Range("C$r").Value2
    $ws2.Range("D$r").Value = $ws1.Range("D$r").Value2
    $ws2.Range("E$r").Value = $ws1.Range("E$r").Value2
}

# --- Update selections ------------------------------------------------------
# Sheet1 is no longer the active tab; its whole used range becomes selected.
$ws1.Range("A1:XFD5").Select()

# Sheet2 becomes the active tab, with the cursor parked at G10.
$ws2.Select()
$ws2.Range("G10").Select()
